$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.141.47"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "2.431.63"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.03"
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.45"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.35"
$ws.Range("E10").Value = "  +3.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0800"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("E12").Value = "  +2.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.70"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.96"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Value = "2.806.33"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "2.407.16"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.832"
$ws.Range("E17").Value = "  +2.56%  "
$ws.Range("D18").Value = "44.128.35"
$ws.Range("E18").Value = "  +2.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.25"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").Value = "0.0₃0908"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.38"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.82"
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.19"
$ws.Range("E27").Value = "  +1.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.35"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").Value = "  +3.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.81"
$ws.Range("E30").Value = "  +3.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("E31").Value = "  +16.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.57"
$ws.Range("E32").Value = "  +7.64%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0761"
$ws.Range("E35").Value = "  +3.48%  "
$ws.Range("E36").Value = "  +3.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "130.09"
$ws.Range("E37").Value = "  +21.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.49"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.93"
$ws.Range("E39").Value = "  +4.75%  "
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.46"
$ws.Range("E42").Value = "  -4.39%  "
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("D44").Value = "1.950.86"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.88"
$ws.Range("E46").Value = "  +4.25%  "
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D48").Value = "2.668.93"
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.62"
$ws.Range("E49").Value = "  +7.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.12"
$ws.Range("E51").Value = "  +1.25%  "
